# This script applies a re-ordering / correction of the observation rows
# 19-22 (Hydnellum aurantiacum / Boletopsis leucomelaena group) and rows
# 23-30 (Hepatica nobilis / Clavariadelphus truncatus group) in the
# "Artfynd" sheet, matching the corrected source export:
#   - row contents are permuted among themselves (ids/measurements that
#     had been attached to the wrong row are moved to the right one)
#   - the Ost/Nord (Q/R) coordinates are stored rounded to whole metres
#     instead of the long floating point values
#   - Starttid/Sluttid (Z/AB) of "00:00" are cleared (no time recorded)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by letter) that carry data for the rows being touched, and
# whether that column should be forced to text so Excel does not
# reinterpret numeric-looking / date-looking strings as real numbers or
# dates (which is how these columns are actually stored in the source).
$columns = @(
    @{ Name = 'A';  Text = $false },
    @{ Name = 'B';  Text = $false },
    @{ Name = 'C';  Text = $true  },
    @{ Name = 'D';  Text = $true  },
    @{ Name = 'E';  Text = $false },
    @{ Name = 'F';  Text = $true  },
    @{ Name = 'G';  Text = $true  },
    @{ Name = 'H';  Text = $true  },
    @{ Name = 'I';  Text = $true  },
    @{ Name = 'J';  Text = $true  },
    @{ Name = 'K';  Text = $true  },
    @{ Name = 'P';  Text = $true  },
    @{ Name = 'Q';  Text = $false },
    @{ Name = 'R';  Text = $false },
    @{ Name = 'S';  Text = $false },
    @{ Name = 'T';  Text = $true  },
    @{ Name = 'U';  Text = $true  },
    @{ Name = 'V';  Text = $true  },
    @{ Name = 'W';  Text = $true  },
    @{ Name = 'Y';  Text = $true  },
    @{ Name = 'Z';  Text = $true  },
    @{ Name = 'AA'; Text = $true  },
    @{ Name = 'AB'; Text = $true  },
    @{ Name = 'AC'; Text = $true  },
    @{ Name = 'AD'; Text = $false },
    @{ Name = 'AE'; Text = $false },
    @{ Name = 'AG'; Text = $false },
    @{ Name = 'AT'; Text = $true  },
    @{ Name = 'AW'; Text = $true  },
    @{ Name = 'AX'; Text = $true  },
    @{ Name = 'AY'; Text = $true  }
)

# Snapshot every value for rows 19-30 before any writes happen, since the
# target rows overlap with the source rows used to populate them.
$snapshot = @{}
foreach ($r in 19..30) {
    $rowData = @{}
    foreach ($col in $columns) {
        $rowData[$col.Name] = $ws.Range("$($col.Name)$r").Value()
    }
    $snapshot[$r] = $rowData
}

# Destination row -> source row (which row's data ends up at the destination).
$mapping = @{
    19 = 22
    20 = 21
    21 = 20
    22 = 19
    23 = 26
    24 = 29
    25 = 23
    26 = 30
    27 = 28
    28 = 24
    29 = 27
    30 = 25
}

foreach ($destRow in 19..30) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]

    foreach ($col in $columns) {
        $name = $col.Name
        $value = $srcData[$name]

        if ($name -eq 'Q' -or $name -eq 'R') {
            if ($null -ne $value -and $value -ne '') {
                $value = [Math]::Round([double]$value)
            }
        }
        elseif ($name -eq 'Z' -or $name -eq 'AB') {
            if ($value -eq '00:00') {
                $value = ''
            }
        }

        $target = $ws.Range("$name$destRow")

        if ($col.Text) {
            # Keep as text even when the content looks numeric/date-like
            # (e.g. "10" or "2023-09-06"), matching the source formatting.
            $target.NumberFormat = "@"
        }

        if ($null -eq $value) {
            $target.Value = ''
        }
        else {
            $target.Value = $value
        }
    }
}
